$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.422.56'
$ws.Range('E2').Value = '  -3.19%  '
$ws.Range('D3').Value = '2.249.42'
$ws.Range('E3').Value = '  -4.09%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.03'
$ws.Range('E5').Value = '  -2.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.633'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '69.72'
$ws.Range('E7').Value = '  -3.65%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.566'
$ws.Range('E9').Value = '  -4.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0992'
$ws.Range('E10').Value = '  -0.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.16'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '35.83'
$ws.Range('E12').Value = '  +10.28%  '
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.77'
$ws.Range('E14').Value = '  -5.51%  '
$ws.Range('D15').Value = '2.578.27'
$ws.Range('E15').Value = '  -4.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.00'
$ws.Range('E16').Value = '  -7.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.858'
$ws.Range('E17').Value = '  -5.12%  '
$ws.Range('D18').Value = '2.246.01'
$ws.Range('E18').Value = '  -4.41%  '
$ws.Range('D19').Value = '42.224.00'
$ws.Range('E19').Value = '  -3.36%  '
$ws.Range('D20').Value = '0.0₃0976'
$ws.Range('E20').Value = '  -4.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.28'
$ws.Range('E21').Value = '  -6.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.44'
$ws.Range('E22').Value = '  -6.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.46'
$ws.Range('E23').Value = '  -6.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.97'
$ws.Range('E24').Value = '  +2.51%  '
$ws.Range('E26').Value = '  -1.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.37'
$ws.Range('E27').Value = '  -5.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.05'
$ws.Range('E28').Value = '  -3.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.18'
$ws.Range('E29').Value = '  -4.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.04'
$ws.Range('E30').Value = '  -3.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.58'
$ws.Range('E31').Value = '  -7.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.121'
$ws.Range('E32').Value = '  -5.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.126'
$ws.Range('E33').Value = '  -6.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0720'
$ws.Range('E34').Value = '  -2.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.36'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.74'
$ws.Range('E36').Value = '  -7.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.63'
$ws.Range('E37').Value = '  -3.14%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '21.56'
$ws.Range('E38').Value = '  +14.09%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.26'
$ws.Range('E39').Value = '  -5.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0269'
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.95'
$ws.Range('E41').Value = '  -6.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '66.39'
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.88'
$ws.Range('E43').Value = '  -6.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.94'
$ws.Range('E44').Value = '  -3.10%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.102'
$ws.Range('E45').Value = '  -4.86%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.189'
$ws.Range('E46').Value = '  -4.12%  '
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('B48').Value = 'BitTorrent-New'
$ws.Range('C48').Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range('D48').Value = '0.0₃0156'
$ws.Range('E48').Value = '  +19.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.45'
$ws.Range('E49').Value = '  +10.03%  '
$ws.Range('E50').Value = '  -3.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.35'
$ws.Range('E51').Value = '  -4.59%  '
